$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H26").Value = 8726.666999999999
$ws.Range("J26").Value = 13015
$ws.Range("L26").Value = 13015
$ws.Range("N26").Value = -13703

$ws.Range("H32").Value = 2000
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").Value = ""

$ws.Range("H33").Value = 247.41176
$ws.Range("I33").Value = 221.64285
$ws.Range("K33").Value = 221.64285
$ws.Range("M33").Value = 7.35714999999999

$ws.Range("H38").Value = 1285.7142
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 1285.7142
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 3857.1426
$ws.Range("M38").Value = ""
$ws.Range("N38").Value = -4601.142599999999

$ws.Range("H137").Value = 2703.8572
$ws.Range("J137").Value = 3325.0833
$ws.Range("L137").Value = 9975.249899999999
$ws.Range("N137").Value = -15075.2499

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 19487.5
$ws.Range("I46").Value = 19000
$ws.Range("J46").Value = 19650
$ws.Range("K46").Value = 19000
$ws.Range("L46").Value = 19650
$ws.Range("M46").Value = -18681
$ws.Range("N46").Value = -20288

$ws.Range("H63").Value = 1998.875
$ws.Range("I63").Value = 1070.2858
$ws.Range("J63").Value = 8499
$ws.Range("K63").Value = 1070.2858
$ws.Range("L63").Value = 8499
$ws.Range("M63").Value = -384.2858000000001
$ws.Range("N63").Value = -9871

$ws.Range("H66").Value = 1998.875
$ws.Range("I66").Value = 1070.2858
$ws.Range("J66").Value = 8499
$ws.Range("K66").Value = 5351.429
$ws.Range("L66").Value = 42495
$ws.Range("M66").Value = -1919.429
$ws.Range("N66").Value = -49359

$ws.Range("H74").Value = 813.7778
$ws.Range("I74").Value = 813.7778
$ws.Range("K74").Value = 813.7778
$ws.Range("M74").Value = 60.22220000000004

$ws.Range("H77").Value = 813.7778
$ws.Range("I77").Value = 813.7778
$ws.Range("K77").Value = 4068.889
$ws.Range("M77").Value = 299.1110000000003

$ws.Range("H132").Value = 4710.25
$ws.Range("I132").Value = 4547
$ws.Range("J132").Value = 5200
$ws.Range("K132").Value = 13641
$ws.Range("L132").Value = 15600
$ws.Range("M132").Value = -11111
$ws.Range("N132").Value = -20660

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2747.6
$ws.Range("I20").Value = 1586.8182
$ws.Range("K20").Value = 1586.8182
$ws.Range("M20").Value = -1339.8182

$ws.Range("H22").Value = 333
$ws.Range("I22").Value = 333
$ws.Range("K22").Value = 333
$ws.Range("M22").Value = -160

$ws.Range("H82").Value = 11300
$ws.Range("I82").Value = 11300
$ws.Range("K82").Value = 11300
$ws.Range("M82").Value = -10917

$ws.Range("H85").Value = 11300
$ws.Range("I85").Value = 11300
$ws.Range("K85").Value = 11300
$ws.Range("M85").Value = -9974

$ws.Range("H86").Value = 2610.8667
$ws.Range("I86").Value = 2733.1428
$ws.Range("J86").Value = 899
$ws.Range("K86").Value = 2733.1428
$ws.Range("L86").Value = 899
$ws.Range("M86").Value = -1610.1428
$ws.Range("N86").Value = -3145

$ws.Range("H89").Value = 2610.8667
$ws.Range("I89").Value = 2733.1428
$ws.Range("J89").Value = 899
$ws.Range("K89").Value = 13665.714
$ws.Range("L89").Value = 4495
$ws.Range("M89").Value = -8049.714
$ws.Range("N89").Value = -15727

$ws.Range("H105").Value = 5325.3335
$ws.Range("I105").Value = 4704.25
$ws.Range("K105").Value = 4704.25
$ws.Range("M105").Value = -2957.25

$ws.Range("H134").Value = 2774.75
$ws.Range("I134").Value = 2034
$ws.Range("K134").Value = 6102
$ws.Range("M134").Value = -3567

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 25499.666
$ws.Range("J43").Value = 25499.666
$ws.Range("L43").Value = 25499.666
$ws.Range("N43").Value = -25867.666

$ws.Range("H59").Value = 44999
$ws.Range("J59").Value = 44999
$ws.Range("L59").Value = 44999
$ws.Range("N59").Value = -47289

$ws.Range("H70").Value = 24999.334
$ws.Range("J70").Value = 24999.334
$ws.Range("L70").Value = 24999.334
$ws.Range("N70").Value = -25629.334

$ws.Range("H73").Value = 24999.334
$ws.Range("J73").Value = 24999.334
$ws.Range("L73").Value = 24999.334
$ws.Range("N73").Value = -27183.334

$ws.Range("H86").Value = 13944694
$ws.Range("I86").Value = 23237822
$ws.Range("K86").Value = 23237822
$ws.Range("M86").Value = -23236699

$ws.Range("H89").Value = 13944694
$ws.Range("I89").Value = 23237822
$ws.Range("K89").Value = 116189110
$ws.Range("M89").Value = -116183494

$ws.Range("H101").Value = 25499.666
$ws.Range("J101").Value = 25499.666
$ws.Range("L101").Value = 25499.666
$ws.Range("N101").Value = -31989.666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 9000
$ws.Range("J12").Value = 9000
$ws.Range("L12").Value = 9000
$ws.Range("N12").Value = -9280

$ws.Range("H70").Value = 8340486.5
$ws.Range("I70").Value = 16674974
$ws.Range("J70").Value = 5999
$ws.Range("K70").Value = 16674974
$ws.Range("L70").Value = 5999
$ws.Range("M70").Value = -16674704
$ws.Range("N70").Value = -6539

$ws.Range("H73").Value = 8340486.5
$ws.Range("I73").Value = 16674974
$ws.Range("J73").Value = 5999
$ws.Range("K73").Value = 16674974
$ws.Range("L73").Value = 5999
$ws.Range("M73").Value = -16674038
$ws.Range("N73").Value = -7871

$ws.Range("H97").Value = 1724.8
$ws.Range("I97").Value = 1031
$ws.Range("K97").Value = 1031
$ws.Range("M97").Value = -535

$ws.Range("H122").Value = 2618.3572
$ws.Range("I122").Value = 2618.3572
$ws.Range("K122").Value = 7855.071599999999
$ws.Range("M122").Value = -5405.071599999999

$ws.Range("H126").Value = 76927450
$ws.Range("I126").Value = 111115176
$ws.Range("J126").Value = 5058
$ws.Range("K126").Value = 333345528
$ws.Range("L126").Value = 15174
$ws.Range("M126").Value = -333343058
$ws.Range("N126").Value = -20114

$ws.Range("H138").Value = 57749.5
$ws.Range("J138").Value = 57749.5
$ws.Range("L138").Value = 57749.5
$ws.Range("N138").Value = -68029.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1059.6666
$ws.Range("I16").Value = 1092
$ws.Range("J16").Value = 995
$ws.Range("K16").Value = 1092
$ws.Range("L16").Value = 995
$ws.Range("M16").Value = -922
$ws.Range("N16").Value = -1335

$ws.Range("H55").Value = 895
$ws.Range("J55").Value = 973.1579
$ws.Range("L55").Value = 973.1579
$ws.Range("N55").Value = -1319.1579

$ws.Range("H132").Value = 3255
$ws.Range("I132").Value = 3408.4
$ws.Range("K132").Value = 10225.2
$ws.Range("M132").Value = -7695.200000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 22000

$ws.Range("H73").Value = 22000

$ws.Range("H92").Value = 49999.5
$ws.Range("J92").Value = 49999.5
$ws.Range("L92").Value = 49999.5
$ws.Range("N92").Value = -54991.5

$ws.Range("H132").Value = 113029.89
$ws.Range("I132").Value = 113029.89
$ws.Range("K132").Value = 339089.67
$ws.Range("M132").Value = -336559.67
